# drawingML export: zero rotation shouldn't alter the shape position
#
# The document contains a legacy VML drawing (w:pict) with a v:group that
# holds a couple of v:rect textboxes. The target edit adds three more VML
# shapes to that same v:group, right after the "Hardware" textbox's
# closing </v:rect> and before the group's closing </v:group>:
#   1. a <v:shapetype> describing an elbow connector
#   2. a <v:rect> textbox containing "MSCAVSS"
#   3. a <v:shape> elbow connector between the two new boxes
#   4. a <v:rect> textbox containing "App-V Agent"
#
# Because this content lives inside a w:pict (pure VML), it is not exposed
# through the normal Shapes/InlineShapes/Find object model (Shapes.Count
# is 0 for this document) - so we locate the host paragraph, splice the
# new VML markup into its raw OOXML as text, and push the result back in
# with Range.InsertXML (the supported way to replace a range's raw
# content).

$d = $word.ActiveDocument

# The exact VML/OOXML to splice in, verbatim from the target markup.
$insertion = @'
<v:shapetype id="_x0000_t34" coordsize="21600,21600" o:spt="34" o:oned="t" adj="10800" path="m,l@0,0@0,21600,21600,21600e" filled="f">
              <v:stroke joinstyle="miter"/>
              <v:formulas>
                <v:f eqn="val #0"/>
              </v:formulas>
              <v:path arrowok="t" fillok="f" o:connecttype="none"/>
              <v:handles>
                <v:h position="#0,center"/>
              </v:handles>
              <o:lock v:ext="edit" shapetype="t"/>
            </v:shapetype>
            <v:rect id="_x0000_s1037" style="position:absolute;left:6823;top:258;width:2095;height:332">
              <v:textbox>
                <w:txbxContent>
                  <w:p>
                    <w:pPr>
                      <w:jc w:val="center"/>
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                    </w:pPr>
                    <w:r>
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                      <w:t>MSCAVSS</w:t>
                    </w:r>
                  </w:p>
                </w:txbxContent>
              </v:textbox>
            </v:rect>
            <v:shape id="_x0000_s1047" type="#_x0000_t34" style="position:absolute;left:4953;top:424;width:1870;height:917;rotation:180;flip:y" o:connectortype="elbow" adj=",254405,-84343">
              <v:stroke startarrow="block" endarrow="block"/>
            </v:shape>
            <v:rect id="_x0000_s1049" style="position:absolute;left:2859;top:1161;width:2094;height:358">
              <v:textbox>
                <w:txbxContent>
                  <w:p>
                    <w:pPr>
                      <w:jc w:val="center"/>
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                    </w:pPr>
                    <w:r>
                      <w:rPr>
                        <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                        <w:sz w:val="16"/>
                        <w:szCs w:val="16"/>
                      </w:rPr>
                      <w:t>App-V Agent</w:t>
                    </w:r>
                  </w:p>
                </w:txbxContent>
              </v:textbox>
            </v:rect>
'@

# Helper: given raw xml text and the index of a tag's opening "<", return
# the index just past the matching close tag - accounting for any nested
# occurrences of the same element name (e.g. nested <w:p> inside a
# txbxContent).
function Get-MatchingCloseEnd {
    param(
        [string]$Xml,
        [int]$OpenIndex,
        [string]$OpenTag,
        [string]$CloseTag
    )

    $pos = $OpenIndex
    $depth = 0
    while ($true) {
        $nextOpen = $Xml.IndexOf($OpenTag, $pos)
        $nextClose = $Xml.IndexOf($CloseTag, $pos)
        if ($nextClose -eq -1) {
            return -1
        }
        if ($nextOpen -ne -1 -and $nextOpen -lt $nextClose) {
            $afterChar = $Xml.Substring($nextOpen + $OpenTag.Length, 1)
            if ($afterChar -eq ">" -or $afterChar -eq " ") {
                $depth = $depth + 1
            }
            $pos = $nextOpen + $OpenTag.Length
        }
        else {
            $depth = $depth - 1
            $pos = $nextClose + $CloseTag.Length
            if ($depth -eq 0) {
                return $pos
            }
        }
    }
}

# Find the paragraph that hosts the VML drawing with the "Hardware" textbox.
$targetParagraph = $null
$targetXml = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $candidate = $p.Range.WordOpenXML
    if ($candidate.Contains("<w:pict>") -and $candidate.Contains(">Hardware<")) {
        $targetParagraph = $p
        $targetXml = $candidate
    }
}

if ($targetParagraph -eq $null) {
    throw "Could not locate the paragraph containing the VML 'Hardware' textbox"
}

# Pull the single <w:p>...</w:p> body paragraph out of the mini-package
# that Range.WordOpenXML hands back.
$bodyTag = "<w:body>"
$bodyStart = $targetXml.IndexOf($bodyTag) + $bodyTag.Length
$pStart = $targetXml.IndexOf("<w:p", $bodyStart)
$pEnd = Get-MatchingCloseEnd $targetXml $pStart "<w:p" "</w:p>"
if ($pEnd -eq -1) {
    throw "Could not find the end of the host paragraph"
}
$paragraphXml = $targetXml.Substring($pStart, $pEnd - $pStart)

# Range.WordOpenXML round-trips the paragraph through a fresh mini-package
# and stamps synthetic w14:paraId/w14:textId attributes along the way;
# strip those back out so we don't introduce attributes that weren't part
# of the original markup (and aren't part of the target edit either).
$paragraphXml = $paragraphXml -replace ' w14:paraId="[0-9A-Fa-f]+"', ''
$paragraphXml = $paragraphXml -replace ' w14:textId="[0-9A-Fa-f]+"', ''

# Splice the new VML shapes in right after the "Hardware" v:rect's
# </v:rect>, i.e. right before the </v:group> that closes the drawing.
$closeGroupTag = "</v:group>"
$hardwareIdx = $paragraphXml.IndexOf(">Hardware<")
$groupCloseIdx = $paragraphXml.IndexOf($closeGroupTag, $hardwareIdx)
if ($groupCloseIdx -eq -1) {
    throw "Could not find the </v:group> closing the drawing"
}

$updatedParagraphXml = $paragraphXml.Substring(0, $groupCloseIdx) + $insertion + "`r`n          " + $paragraphXml.Substring($groupCloseIdx)

# Declare the namespaces the fragment relies on directly on the root
# element so InsertXML can resolve v:/o:/w: prefixes on their own.
$updatedParagraphXml = $updatedParagraphXml -replace '^<w:p ', '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" '

$targetParagraph.Range.InsertXML($updatedParagraphXml)

"Inserted MSCAVSS/App-V Agent VML shapes into the drawing"
